$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '28.997.69'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -1.30%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.913.59'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  -0.37%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '324.89'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  -0.33%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4591'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.85%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3826'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.23%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07723'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -1.40%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.9817'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '22.15'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -2.50%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.900.02'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '6.951'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.08%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.673'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.43%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.07026'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E16").Value = '  -0.39%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '83.99'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -3.44%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000009466'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -3.99%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '16.73'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.61%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.43%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '28.983.26'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.44%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.323'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.02%  '
$ws.Range("E23").Value = '  -1.56%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.089'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.61%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '158.83'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.78%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '19.07'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.84%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '5.680'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.46%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '117.50'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.03%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.859'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.09302'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.85%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.8677'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.46%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.085'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.74%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.256'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.69%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '3.025'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.21%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.05736'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.32%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.156'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("E38").Value = '  -1.89%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.5513'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.82%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '7.421'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -3.47%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.1757'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.40%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.871'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +5.30%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '9.353'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.73%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.5192'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.01%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '11.29'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.86%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.06870'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.14%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.000002632'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -7.60%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.060'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.74%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.784'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.93%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '110.63'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.13%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.36%  '
